$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,4,5,7,8,9,11,12,13)  # B,C,D,E,G,H,I,K,L,M
$data = @(
    @(0.596238626253438, 0.05275607247739345, 0.07707292167084745, 0.08133836603712652, 2.118018728967002, 1.754318280332768, 1.810319534582341, 0.4212422880382007, 0.2267735543746383, 0.1815961741146062),
    @(0.5709050036153656, 0.04803073841155481, 0.07012574728479137, 0.08141792776114087, 2.094241755356848, 1.747844621939265, 1.800593244768102, 0.3920167350184727, 0.2239746217362182, 0.1764762417096897),
    @(0.555746512494494, 0.04510287436224303, 0.06589753627913808, 0.08148659240934286, 2.080463671151008, 1.744431898839963, 1.795207060022165, 0.3743417770515407, 0.2223646465878204, 0.1734433447835784),
    @(0.5496691913925247, 0.04390301604622948, 0.06418385353782696, 0.08151956608147515, 2.075055173692476, 1.743182527198641, 1.793159356247955, 0.3672069311591883, 0.2217358993454468, 0.1722353099306879),
    @(0.5486660956487128, 0.04370337200859353, 0.0638998608739314, 0.08152534315119375, 2.074169541624542, 1.742983606342008, 1.792828225053483, 0.3660262942128298, 0.221633147866477, 0.1720364027651016),
    @(0.5556641469156318, 0.04508671999292346, 0.06587438717495786, 0.08148701687680671, 2.080389895813084, 1.744414477117999, 1.795178848055109, 0.374245279278, 0.2223560563850313, 0.1734269397981336),
    @(0.587421414807352, 0.05113223256381616, 0.07466974344303878, 0.08136169255896419, 2.109649833617965, 1.751969480918007, 1.806844265710957, 0.4111093625356546, 0.2257859612031297, 0.1798078562333423),
    @(0.6528387605580406, 0.06278087174425195, 0.09221731823673451, 0.08127274076810487, 2.173560903178782, 1.771247838970964, 1.834374982626116, 0.485542147665825, 0.2333729744536441, 0.1931987115020384),
    @(0.7028166573529404, 0.07121885461336319, 0.1052982582023105, 0.08130249949973134, 2.224530028537288, 1.788139543213561, 1.857452807817566, 0.5415445712513076, 0.2394722177147202, 0.2035722857155946),
    @(0.7259695090946536, 0.07503265881992149, 0.1112914503664229, 0.08133659242415447, 2.248595925655394, 1.796418121398005, 1.868573690216465, 0.5673101598346193, 0.2423610435771479, 0.2084078620727325),
    @(0.7347968720181086, 0.07647338489518063, 0.1135671224978978, 0.08135244942837261, 2.257836004929629, 1.799638568147742, 1.87287459428461, 0.5771086791311859, 0.2434713860890412, 0.2102557164773629),
    @(0.7328930810577106, 0.0761632523654896, 0.1130767400617714, 0.08134890342564383, 2.255840341668772, 1.798941182639197, 1.871944327599678, 0.5749965406161266, 0.2432315246106214, 0.2098570046494714),
    @(0.7266945424656512, 0.07515125754170526, 0.1114785471119433, 0.08133783798395378, 2.249353568430536, 1.796681355092545, 1.86892573037575, 0.5681154553049623, 0.2424520634308323, 0.2085595511570091),
    @(0.722905551594522, 0.07453093052106396, 0.1105004142363555, 0.08133144356667721, 2.24539676159776, 1.795308285260063, 1.867088433581401, 0.5639060153046387, 0.241976756960625, 0.2077670006044343),
    @(0.7013119530327003, 0.07096912293005175, 0.104907450929943, 0.08130068429708004, 2.222974997642694, 1.787610484610809, 1.856738572590288, 0.5398665589606821, 0.2392857232438672, 0.2032586121347322),
    @(0.6881718234071457, 0.06877779853878963, 0.1014872997693885, 0.08128707274737224, 2.209445543389705, 1.783040420728668, 1.850548838748935, 0.5251933391847103, 0.2376641073980892, 0.2005226946576428),
    @(0.6806532904949734, 0.06751507985462979, 0.09952412644696551, 0.081281179068867, 2.201746548716187, 1.780467791070009, 1.847047267846108, 0.5167809702622321, 0.2367421518515016, 0.1989600430909562),
    @(0.678114403541457, 0.06708714324862797, 0.09886011629717473, 0.08127951624276442, 2.199154010446591, 1.779606350988786, 1.845871756855374, 0.5139373726658221, 0.2364318414062581, 0.1984328424054453),
    @(0.6895665442506527, 0.06901130907877473, 0.1018509655874027, 0.08128832148750043, 2.210877206059308, 1.783521120851276, 1.851201680634531, 0.5267525053893394, 0.2378356182335324, 0.2008128020693505),
    @(0.728513579762847, 0.07544859880984234, 0.1119478069634425, 0.081341008270714, 2.251255445728674, 1.797342799613688, 1.869809930975961, 0.5701354666577174, 0.2426805649700157, 0.2089401910443769),
    @(0.7543166685591416, 0.07963549518562729, 0.1185827251215699, 0.08139261388024721, 2.278384473187771, 1.806874616065301, 1.882494208212307, 0.5987315002184062, 0.2459426395250688, 0.2143493877749876),
    @(0.7405132031938422, 0.07740269907165498, 0.1150382282504978, 0.08136350277209203, 2.26383742433751, 1.801741673057705, 1.875676499899726, 0.5834470530708415, 0.2441928667168582, 0.2114534908265497),
    @(0.6889358791112556, 0.06890574799398053, 0.1016865425679612, 0.08128775091438101, 2.210229704464155, 1.783303625735698, 1.850906353296466, 0.5260475332720773, 0.2377580459981488, 0.2006816125353481),
    @(0.6348052702464599, 0.05965113824927926, 0.08743751743706696, 0.08128006729614867, 2.155568692323584, 1.76555400515133, 1.826427489813568, 0.46517566731643, 0.23122829483777, 0.1894821360281931)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $ws.Cells.Item($r + 2, $cols[$c]).Value2 = $data[$r][$c]
    }
}